$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 24446.428
$ws.Range("I33").Value = 25663.25
$ws.Range("J33").Value = 110
$ws.Range("K33").Value = 25663.25
$ws.Range("L33").Value = 110
$ws.Range("M33").Value = -25434.25
$ws.Range("N33").Value = -568

$ws.Range("H74").Value = 5362.5
$ws.Range("I74").Value = 4928.7144
$ws.Range("J74").Value = 8399
$ws.Range("K74").Value = 4928.7144
$ws.Range("L74").Value = 8399
$ws.Range("M74").Value = -3992.7144
$ws.Range("N74").Value = -10271

$ws.Range("H77").Value = 5362.5
$ws.Range("I77").Value = 4928.7144
$ws.Range("J77").Value = 8399
$ws.Range("K77").Value = 24643.572
$ws.Range("L77").Value = 41995
$ws.Range("M77").Value = -19963.572
$ws.Range("N77").Value = -51355

$ws.Range("H111").Value = 1257.2
$ws.Range("I111").Value = 1257.2
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3771.6
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -704.6000000000004

$ws.Range("H112").Value = 2622.125
$ws.Range("I112").Value = 1918.3334
$ws.Range("J112").Value = 3044.4
$ws.Range("K112").Value = 5755.0002
$ws.Range("L112").Value = 9133.200000000001
$ws.Range("M112").Value = -4647.0002
$ws.Range("N112").Value = -11349.2

$ws.Range("H137").Value = 23685806
$ws.Range("I137").Value = 75001096
$ws.Range("J137").Value = 1825.3846
$ws.Range("K137").Value = 225003288
$ws.Range("L137").Value = 5476.1538
$ws.Range("M137").Value = -225000738
$ws.Range("N137").Value = -10576.1538

$ws.Range("H138").Value = 4304.273
$ws.Range("I138").Value = 5269.6
$ws.Range("J138").Value = 3499.8333
$ws.Range("K138").Value = 15808.8
$ws.Range("L138").Value = 10499.4999
$ws.Range("M138").Value = -10668.8
$ws.Range("N138").Value = -20779.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 7995.5
$ws.Range("I46").Value = 8498
$ws.Range("J46").Value = 7493
$ws.Range("K46").Value = 8498
$ws.Range("L46").Value = 7493
$ws.Range("M46").Value = -8179
$ws.Range("N46").Value = -8131

$ws.Range("H74").Value = 621594.3
$ws.Range("I74").Value = 4312.4736
$ws.Range("J74").Value = 3972553
$ws.Range("K74").Value = 4312.4736
$ws.Range("L74").Value = 3972553
$ws.Range("M74").Value = -3438.4736
$ws.Range("N74").Value = -3974301

$ws.Range("H77").Value = 621594.3
$ws.Range("I77").Value = 4312.4736
$ws.Range("J77").Value = 3972553
$ws.Range("K77").Value = 21562.368
$ws.Range("L77").Value = 19862765
$ws.Range("M77").Value = -17194.368
$ws.Range("N77").Value = -19871501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2227.36
$ws.Range("I86").Value = 1861.7
$ws.Range("J86").Value = 3690
$ws.Range("K86").Value = 1861.7
$ws.Range("L86").Value = 3690
$ws.Range("M86").Value = -738.7

$ws.Range("H89").Value = 2227.36
$ws.Range("I89").Value = 1861.7
$ws.Range("J89").Value = 3690
$ws.Range("K89").Value = 9308.5
$ws.Range("L89").Value = 18450
$ws.Range("M89").Value = -3692.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N16").ClearContents()
$ws.Range("H16").Value = 1057
$ws.Range("I16").Value = 1057
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1057
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -770

$ws.Range("H31").Value = 6387.3237
$ws.Range("I31").Value = 3910.6667
$ws.Range("J31").Value = 7278.92
$ws.Range("K31").Value = 3910.6667
$ws.Range("L31").Value = 7278.92
$ws.Range("M31").Value = -3615.6667
$ws.Range("N31").Value = -7868.92

$ws.Range("H34").Value = 6387.3237
$ws.Range("I34").Value = 3910.6667
$ws.Range("J34").Value = 7278.92
$ws.Range("K34").Value = 3910.6667
$ws.Range("L34").Value = 7278.92
$ws.Range("M34").Value = -3708.6667
$ws.Range("N34").Value = -7682.92

$ws.Range("N113").ClearContents()
$ws.Range("H113").Value = 1057
$ws.Range("I113").Value = 1057
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1057
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1113

$ws.Range("H124").Value = 78982
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 78982
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 78982
$ws.Range("N124").Value = -83892

$ws.Range("H132").Value = 3240.2104
$ws.Range("I132").Value = 3337.5
$ws.Range("J132").Value = 1489
$ws.Range("K132").Value = 10012.5
$ws.Range("L132").Value = 4467
$ws.Range("M132").Value = -7482.5

$ws.Range("N133").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2054.4
$ws.Range("I5").Value = 1113.4286
$ws.Range("J5").Value = 4250
$ws.Range("K5").Value = 3340.2858
$ws.Range("L5").Value = 12750
$ws.Range("M5").Value = -3228.2858
$ws.Range("N5").Value = -12974

$ws.Range("H46").Value = 9093317
$ws.Range("I46").Value = 16667917
$ws.Range("J46").Value = 3797.8
$ws.Range("K46").Value = 50003751
$ws.Range("L46").Value = 11393.4
$ws.Range("M46").Value = -50003660
$ws.Range("N46").Value = -11575.4

$ws.Range("H68").Value = 3268.52
$ws.Range("I68").Value = 1578.4286
$ws.Range("J68").Value = 3925.7778
$ws.Range("K68").Value = 4735.2858
$ws.Range("L68").Value = 11777.3334
$ws.Range("M68").Value = -3924.2858
$ws.Range("N68").Value = -13399.3334

$ws.Range("H71").Value = 3268.52
$ws.Range("I71").Value = 1578.4286
$ws.Range("J71").Value = 3925.7778
$ws.Range("K71").Value = 14205.8574
$ws.Range("L71").Value = 35332.00019999999
$ws.Range("M71").Value = -10149.8574
$ws.Range("N71").Value = -43444.00019999999

$ws.Range("H94").Value = 10673.833
$ws.Range("I94").Value = 11010.75
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 33032.25
$ws.Range("L94").Value = 30000
$ws.Range("M94").Value = -32356.25

$ws.Range("H97").Value = 54030.707
$ws.Range("I97").Value = 72722.91
$ws.Range("J97").Value = 19761.666
$ws.Range("K97").Value = 218168.73
$ws.Range("L97").Value = 59284.99800000001
$ws.Range("M97").Value = -217672.73
$ws.Range("N97").Value = -60276.99800000001

$ws.Range("H103").Value = 1341.5714
$ws.Range("I103").Value = 447.5
$ws.Range("J103").Value = 1699.2
$ws.Range("K103").Value = 1342.5
$ws.Range("L103").Value = 5097.6
$ws.Range("M103").Value = -463.5
$ws.Range("N103").Value = -6855.6

$ws.Range("H106").Value = 45000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 45000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 135000
$ws.Range("N106").Value = -136892

$ws.Range("H121").Value = 1620080.8
$ws.Range("I121").Value = 334333
$ws.Range("J121").Value = 2005805.1
$ws.Range("K121").Value = 1002999
$ws.Range("L121").Value = 6017415.300000001
$ws.Range("M121").Value = -1001689
$ws.Range("N121").Value = -6020035.300000001

$ws.Range("H122").Value = 30318.234
$ws.Range("I122").Value = 506.16666
$ws.Range("J122").Value = 46579.363
$ws.Range("K122").Value = 4555.49994
$ws.Range("L122").Value = 419214.267
$ws.Range("M122").Value = -2105.49994
$ws.Range("N122").Value = -424114.267

$ws.Range("H135").Value = 2054.4
$ws.Range("I135").Value = 1113.4286
$ws.Range("J135").Value = 4250
$ws.Range("K135").Value = 10020.8574
$ws.Range("L135").Value = 38250
$ws.Range("M135").Value = -7485.857399999999
$ws.Range("N135").Value = -43320

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 50057
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 50057
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 50057
$ws.Range("N63").Value = -51429

$ws.Range("H66").Value = 50057
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 50057
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 150171
$ws.Range("N66").Value = -157035

$ws.Range("H111").Value = 51999.332
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 51999.332
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 51999.332
$ws.Range("N111").Value = -58133.332

$ws.Range("H122").Value = 3221.8333
$ws.Range("I122").Value = 4123.636
$ws.Range("J122").Value = 1804.7142
$ws.Range("K122").Value = 12370.908
$ws.Range("L122").Value = 5414.142599999999
$ws.Range("M122").Value = -9920.908000000001
$ws.Range("N122").Value = -10314.1426

$ws.Range("H126").Value = 2249.5
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2249.5
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 6748.5
$ws.Range("N126").Value = -11688.5

$ws.Range("N132").ClearContents()
$ws.Range("H132").Value = 19333
$ws.Range("I132").Value = 19333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 57999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -55469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3326.7778
$ws.Range("I7").Value = 3349.1428
$ws.Range("J7").Value = 3248.5
$ws.Range("K7").Value = 3349.1428
$ws.Range("L7").Value = 3248.5
$ws.Range("M7").Value = -3237.1428
$ws.Range("N7").Value = -3472.5

$ws.Range("H40").Value = 3307
$ws.Range("I40").Value = 2995.1333
$ws.Range("J40").Value = 4866.3335
$ws.Range("K40").Value = 2995.1333
$ws.Range("L40").Value = 4866.3335
$ws.Range("M40").Value = -2859.1333
$ws.Range("N40").Value = -5138.3335

$ws.Range("H125").Value = 99000.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 99000.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 99000.5
$ws.Range("N125").Value = -108840.5

$ws.Range("H126").Value = 3326.7778
$ws.Range("I126").Value = 3349.1428
$ws.Range("J126").Value = 3248.5
$ws.Range("K126").Value = 10047.4284
$ws.Range("L126").Value = 9745.5
$ws.Range("M126").Value = -7577.428400000001
$ws.Range("N126").Value = -14685.5

$ws.Range("H136").Value = 76926820
$ws.Range("I136").Value = 3324.5
$ws.Range("J136").Value = 200004400
$ws.Range("K136").Value = 9973.5
$ws.Range("L136").Value = 600013200
$ws.Range("M136").Value = -7423.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 53352.332
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 53352.332
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 53352.332
$ws.Range("N64").Value = -53848.332

$ws.Range("H67").Value = 53352.332
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 53352.332
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 53352.332
$ws.Range("N67").Value = -55068.332

$ws.Range("H81").Value = 29931.666
$ws.Range("I81").Value = 77795
$ws.Range("J81").Value = 6000
$ws.Range("K81").Value = 155590
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = -154529
$ws.Range("N81").Value = -14122

$ws.Range("H84").Value = 29931.666
$ws.Range("I84").Value = 77795
$ws.Range("J84").Value = 6000
$ws.Range("K84").Value = 777950
$ws.Range("L84").Value = 60000
$ws.Range("M84").Value = -772646
$ws.Range("N84").Value = -70608

$ws.Range("N93").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
